$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds plain "YYYY-MM-DD" text in this sheet (ticket log), but Excel
# would otherwise parse such strings into date serials on assignment. Force the
# new rows in column A to Text format first so the literal strings are preserved,
# matching the rest of the sheet (which stores everything as text).
$ws.Range("A60:A67").NumberFormat = "@"

# Row 60
$ws.Cells.Item(60, 1).Value = "2024-05-17"
$ws.Cells.Item(60, 2).Value = "10:56:37"
$ws.Cells.Item(60, 3).Value = "-"
$ws.Cells.Item(60, 4).Value = "Cámara no detecta busbar"
$ws.Cells.Item(60, 5).Value = "-"
$ws.Cells.Item(60, 6).Value = "-"
$ws.Cells.Item(60, 7).Value = "-"
$ws.Cells.Item(60, 8).Value = "10:56:39"
$ws.Cells.Item(60, 9).Value = "0:00:02"

# Row 61
$ws.Cells.Item(61, 1).Value = "2024-05-17"
$ws.Cells.Item(61, 2).Value = "11:02:43"
$ws.Cells.Item(61, 3).Value = "-"
$ws.Cells.Item(61, 4).Value = "Cámara no detecta Pcb"
$ws.Cells.Item(61, 5).Value = "-"
$ws.Cells.Item(61, 6).Value = "-"
$ws.Cells.Item(61, 7).Value = "-"
$ws.Cells.Item(61, 8).Value = "11:02:45"
$ws.Cells.Item(61, 9).Value = "0:00:02"

# Row 62
$ws.Cells.Item(62, 1).Value = "2024-05-17"
$ws.Cells.Item(62, 2).Value = "11:14:24"
$ws.Cells.Item(62, 3).Value = "-"
$ws.Cells.Item(62, 4).Value = "Cámara no detecta Pcb"
$ws.Cells.Item(62, 5).Value = "-"
$ws.Cells.Item(62, 6).Value = "-"
$ws.Cells.Item(62, 7).Value = "-"
$ws.Cells.Item(62, 8).Value = "11:14:25"
$ws.Cells.Item(62, 9).Value = "0:00:01"

# Row 63
$ws.Cells.Item(63, 1).Value = "2024-05-17"
$ws.Cells.Item(63, 2).Value = "11:16:20"
$ws.Cells.Item(63, 3).Value = "-"
$ws.Cells.Item(63, 4).Value = "-"
$ws.Cells.Item(63, 5).Value = "Screw K30 no lo detecta puesto"
$ws.Cells.Item(63, 6).Value = "-"
$ws.Cells.Item(63, 7).Value = "-"
$ws.Cells.Item(63, 8).Value = "11:16:22"
$ws.Cells.Item(63, 9).Value = "0:00:02"

# Row 64
$ws.Cells.Item(64, 1).Value = "2024-05-17"
$ws.Cells.Item(64, 2).Value = "11:16:29"
$ws.Cells.Item(64, 3).Value = "-"
$ws.Cells.Item(64, 4).Value = "-"
$ws.Cells.Item(64, 5).Value = "Detección de sealling mal puesto"
$ws.Cells.Item(64, 6).Value = "-"
$ws.Cells.Item(64, 7).Value = "-"

# Row 65
$ws.Cells.Item(65, 1).Value = "2024-05-17"
$ws.Cells.Item(65, 2).Value = "11:58:14"
$ws.Cells.Item(65, 3).Value = "-"
$ws.Cells.Item(65, 4).Value = "-"
$ws.Cells.Item(65, 5).Value = "Detección de sealling mal puesto"
$ws.Cells.Item(65, 6).Value = "-"
$ws.Cells.Item(65, 7).Value = "-"
$ws.Cells.Item(65, 8).Value = "11:58:18"
$ws.Cells.Item(65, 9).Value = "0:00:04"

# Row 66
$ws.Cells.Item(66, 1).Value = "2024-05-17"
$ws.Cells.Item(66, 2).Value = "11:58:31"
$ws.Cells.Item(66, 3).Value = "-"
$ws.Cells.Item(66, 4).Value = "-"
$ws.Cells.Item(66, 5).Value = "Atasco tuerca"
$ws.Cells.Item(66, 6).Value = "-"
$ws.Cells.Item(66, 7).Value = "-"
$ws.Cells.Item(66, 8).Value = "12:00:10"
$ws.Cells.Item(66, 9).Value = "0:01:39"

# Row 67
$ws.Cells.Item(67, 1).Value = "2024-05-17"
$ws.Cells.Item(67, 2).Value = "12:02:19"
$ws.Cells.Item(67, 3).Value = "-"
$ws.Cells.Item(67, 4).Value = "-"
$ws.Cells.Item(67, 5).Value = "Detección de sealling mal puesto"
$ws.Cells.Item(67, 6).Value = "-"
$ws.Cells.Item(67, 7).Value = "-"
$ws.Cells.Item(67, 8).Value = "12:02:25"
$ws.Cells.Item(67, 9).Value = "0:00:06"

